# "Generate Report for Handback"
#
# The handback transform for the cb64d49f-... file failed because the
# handback file name (sm05robw.oc3) didn't match the handoff file name for
# either locale. Update the localization-status report to reflect that:
#   - flip the row's Status from "Ready for handoff" to
#     "Handback transform failed" everywhere it is shown (Overview summary
#     columns + each locale sheet's Status column), and
#   - record the mismatch message in the Error Detail column of each
#     locale sheet, widening that column so the message is readable.

$wb = $excel.ActiveWorkbook

$statusNew = "Handback transform failed"

$zhError = "Handback file name: sm05robw.oc3 is different with handoff file name: cb64d49f-b0b1-4dd2-a652-7a067b169588.315779d0bd7f5c82242edf0ad640d57c4188d969.zh-cn."
$deError = "Handback file name: sm05robw.oc3 is different with handoff file name: cb64d49f-b0b1-4dd2-a652-7a067b169588.315779d0bd7f5c82242edf0ad640d57c4188d969.de-de."

# --- Overview sheet: zh-cn / de-de status columns for the cb64d49f row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

# --- zh-cn sheet: Status + Error Detail for the cb64d49f row ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusNew
$wsZh.Range("P3").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: Status + Error Detail for the cb64d49f row ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusNew
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = 39.17
